$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 17, pushing the existing rows 17-31 down to 18-32
# (weekly price update: a new week's record is prepended to the series).
$ws.Rows("17:17").Insert()

# Populate the newly inserted row 17 with this week's record.
$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = 44874
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = 300000000
$ws.Range("G17").Value = "Espárragos"
$ws.Range("H17").Value = "Sin especificar"
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 2000
$ws.Range("K17").Value = 1000
$ws.Range("L17").Value = 1100
$ws.Range("M17").Value = 1050
$ws.Range("N17").Value = "$/kilo"
$ws.Range("O17").Value = "Provincia de Diguillín"
$ws.Range("P17").Value = 1050
$ws.Range("Q17").Value = 1
$ws.Range("R17").Value = "Hortaliza"
